$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 237 (shifts existing rows 237-262 down to 238-263)
$ws.Rows("237:237").Insert()

# Populate the newly inserted row 237 with the new weekly price observation
$ws.Cells.Item(237, 1).Value  = 3
$ws.Cells.Item(237, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value  = "Coquimbo"
$ws.Cells.Item(237, 4).Value  = 45106
$ws.Cells.Item(237, 5).Value  = 5
$ws.Cells.Item(237, 6).Value  = 100112026
$ws.Cells.Item(237, 7).Value  = "Haba"
$ws.Cells.Item(237, 8).Value  = "Sin especificar"
$ws.Cells.Item(237, 9).Value  = "Primera"
$ws.Cells.Item(237, 10).Value = 85
$ws.Cells.Item(237, 11).Value = 18500
$ws.Cells.Item(237, 12).Value = 19000
$ws.Cells.Item(237, 13).Value = 18735
$ws.Cells.Item(237, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(237, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(237, 16).Value = 749
$ws.Cells.Item(237, 17).Value = 25
$ws.Cells.Item(237, 18).Value = "Hortaliza"
